$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix bug: chart wasn't displaying for household (ag_comm) row - missing
# "Disaggregates" value in column H (level_lab), matching the other rows.
$ws.Range("H4").Value = "Disaggregates"

# Update the active selection to reflect new interface layout cursor position.
$ws.Range("H5").Select()
